$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4003548.2
$ws.Range("I106").Value = 5003585
$ws.Range("K106").Value = 5003585
$ws.Range("M106").Value = -5002954

$ws.Range("H116").Value = 6750.4546
$ws.Range("I116").Value = 13411.667
$ws.Range("J116").Value = 2138.8462
$ws.Range("K116").Value = 13411.667
$ws.Range("L116").Value = 2138.8462
$ws.Range("M116").Value = -9969.666999999999
$ws.Range("N116").Value = -9022.8462

$ws.Range("H132").Value = 2628.082
$ws.Range("I132").Value = 2513.2407
$ws.Range("J132").Value = 3514
$ws.Range("K132").Value = 7539.722099999999
$ws.Range("L132").Value = 10542
$ws.Range("M132").Value = -5009.722099999999
$ws.Range("N132").Value = -15602

$ws.Range("H138").Value = 1058.93
$ws.Range("I138").Value = 552.1731
$ws.Range("J138").Value = 1607.9166
$ws.Range("K138").Value = 1656.5193
$ws.Range("L138").Value = 4823.7498
$ws.Range("M138").Value = 3483.4807
$ws.Range("N138").Value = -15103.7498


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 644506.2
$ws.Range("I32").Value = 732057.3
$ws.Range("K32").Value = 732057.3
$ws.Range("M32").Value = -731770.3

$ws.Range("H123").Value = 36166.5
$ws.Range("J123").Value = 36166.5
$ws.Range("L123").Value = 36166.5
$ws.Range("N123").Value = -45966.5

$ws.Range("H131").Value = 39735
$ws.Range("J131").Value = 39735
$ws.Range("L131").Value = 39735
$ws.Range("N131").Value = -49815

$ws.Range("H132").Value = 2602.5762
$ws.Range("I132").Value = 2628.8572
$ws.Range("J132").Value = 2564.25
$ws.Range("K132").Value = 7886.571599999999
$ws.Range("L132").Value = 7692.75
$ws.Range("M132").Value = -5356.571599999999
$ws.Range("N132").Value = -12752.75


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4287.477
$ws.Range("I31").Value = 1264.6177
$ws.Range("J31").Value = 7602.871
$ws.Range("K31").Value = 1264.6177
$ws.Range("L31").Value = 7602.871
$ws.Range("M31").Value = -969.6177
$ws.Range("N31").Value = -8192.870999999999

$ws.Range("H34").Value = 4287.477
$ws.Range("I34").Value = 1264.6177
$ws.Range("J34").Value = 7602.871
$ws.Range("K34").Value = 1264.6177
$ws.Range("L34").Value = 7602.871
$ws.Range("M34").Value = -1062.6177
$ws.Range("N34").Value = -8006.871

$ws.Range("H58").Value = 1144.32
$ws.Range("I58").Value = 939.3226
$ws.Range("K58").Value = 939.3226
$ws.Range("M58").Value = -736.3226

$ws.Range("H122").Value = 1861.963
$ws.Range("I122").Value = 1681.625
$ws.Range("J122").Value = 1937.8948
$ws.Range("K122").Value = 5044.875
$ws.Range("L122").Value = 5813.6844
$ws.Range("M122").Value = -2594.875
$ws.Range("N122").Value = -10713.6844

$ws.Range("H132").Value = 4066702.8
$ws.Range("I132").Value = 1400.6207
$ws.Range("J132").Value = 13891183
$ws.Range("K132").Value = 4201.8621
$ws.Range("L132").Value = 41673549
$ws.Range("M132").Value = -1671.8621
$ws.Range("N132").Value = -41678609

$ws.Range("H134").Value = 4036.738
$ws.Range("I134").Value = 3939.257
$ws.Range("K134").Value = 11817.771
$ws.Range("M134").Value = -9282.771000000001

$ws.Range("H136").Value = 1144.32
$ws.Range("I136").Value = 939.3226
$ws.Range("K136").Value = 2817.9678
$ws.Range("M136").Value = -267.9677999999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 457.51614
$ws.Range("I113").Value = 444.66666
$ws.Range("J113").Value = 465.6316
$ws.Range("K113").Value = 1333.99998
$ws.Range("L113").Value = 1396.8948
$ws.Range("M113").Value = 836.0000199999999
$ws.Range("N113").Value = -5736.8948

$ws.Range("H121").Value = 527480
$ws.Range("I121").Value = 3333793
$ws.Range("J121").Value = 1296.2812
$ws.Range("K121").Value = 10001379
$ws.Range("L121").Value = 3888.8436
$ws.Range("M121").Value = -10000069
$ws.Range("N121").Value = -6508.8436


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 12948.8
$ws.Range("J123").Value = 12948.8
$ws.Range("L123").Value = 12948.8
$ws.Range("N123").Value = -17848.8

$ws.Range("H132").Value = 2277.5
$ws.Range("I132").Value = 1860.4043
$ws.Range("J132").Value = 3309.2632
$ws.Range("K132").Value = 5581.2129
$ws.Range("L132").Value = 9927.7896
$ws.Range("M132").Value = -3051.2129
$ws.Range("N132").Value = -14987.7896


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 66006.664
$ws.Range("J4").Value = 66006.664
$ws.Range("L4").Value = 66006.664
$ws.Range("N4").Value = -66232.664

$ws.Range("H28").Value = 66006.664
$ws.Range("J28").Value = 66006.664
$ws.Range("L28").Value = 66006.664
$ws.Range("N28").Value = -66470.664

$ws.Range("H37").Value = 66006.664
$ws.Range("J37").Value = 66006.664
$ws.Range("L37").Value = 66006.664
$ws.Range("N37").Value = -66220.664

$ws.Range("H61").Value = 3112.3704
$ws.Range("I61").Value = 2874.0908
$ws.Range("J61").Value = 4160.8
$ws.Range("K61").Value = 2874.0908
$ws.Range("L61").Value = 4160.8
$ws.Range("M61").Value = -2672.0908
$ws.Range("N61").Value = -4564.8

$ws.Range("H113").Value = 3112.3704
$ws.Range("I113").Value = 2874.0908
$ws.Range("J113").Value = 4160.8
$ws.Range("K113").Value = 2874.0908
$ws.Range("L113").Value = 4160.8
$ws.Range("M113").Value = -704.0907999999999
$ws.Range("N113").Value = -8500.799999999999

$ws.Range("H122").Value = 4473.6665
$ws.Range("I122").Value = 3900
$ws.Range("J122").Value = 4760.5
$ws.Range("K122").Value = 11700
$ws.Range("L122").Value = 14281.5
$ws.Range("M122").Value = -9250
$ws.Range("N122").Value = -19181.5

$ws.Range("H132").Value = 2441.2576
$ws.Range("I132").Value = 2288.2856
$ws.Range("J132").Value = 2882.1765
$ws.Range("K132").Value = 6864.8568
$ws.Range("L132").Value = 8646.529500000001
$ws.Range("M132").Value = -4334.8568
$ws.Range("N132").Value = -13706.5295

$ws.Range("H136").Value = 5378128
$ws.Range("I136").Value = 1600.8422
$ws.Range("J136").Value = 13890962
$ws.Range("K136").Value = 4802.5266
$ws.Range("L136").Value = 41672886
$ws.Range("M136").Value = -2252.5266
$ws.Range("N136").Value = -41677986

$ws.Range("H137").Value = 34333.332
$ws.Range("I137").Value = 23000
$ws.Range("J137").Value = 40000
$ws.Range("K137").Value = 23000
$ws.Range("L137").Value = 40000
$ws.Range("M137").Value = -17900
$ws.Range("N137").Value = -50200


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 2006
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2006
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2006
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = -2286

$ws.Range("H37").Value = 50014.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 50014.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 50014.5
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = -50420.5

$ws.Range("H40").Value = 1800
$ws.Range("I40").Value = 1800
$ws.Range("K40").Value = 1800
$ws.Range("M40").Value = -1651

$ws.Range("H123").Value = 20659.133
$ws.Range("J123").Value = 21977.4
$ws.Range("L123").Value = 21977.4
$ws.Range("N123").Value = -31777.4

$ws.Range("H132").Value = 8104322.5
$ws.Range("I132").Value = 3618.2666
$ws.Range("J132").Value = 13890540
$ws.Range("K132").Value = 10854.7998
$ws.Range("L132").Value = 41671620
$ws.Range("M132").Value = -8324.799800000001
$ws.Range("N132").Value = -41676680

$ws.Range("H136").Value = 2334.65
$ws.Range("I136").Value = 2586.8125
$ws.Range("J136").Value = 1956.4062
$ws.Range("K136").Value = 7760.4375
$ws.Range("L136").Value = 5869.2186
$ws.Range("M136").Value = -5210.4375
$ws.Range("N136").Value = -10969.2186

